$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8490044474601746
$ws.Range("B1").Value = 1.512098670005798
$ws.Range("C1").Value = 6.221972942352295
$ws.Range("D1").Value = 2.935715913772583
$ws.Range("E1").Value = 1.600641012191772
